$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 corresponds to "LeBron James" (ast/blk/pts/stl/trb per game averages).
# Update the values per the diff: B2, C2, D2, and F2 change; E2 stays the same.
$ws.Range("B2").Value = 7.45
$ws.Range("C2").Value = 0.73
$ws.Range("D2").Value = 26.96
$ws.Range("F2").Value = 7.54
